$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 (mirror style/format of existing header cells like H1)
$h1 = $ws.Range("H1")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$headerTarget = $ws.Range("I1:J1")
$headerTarget.Font.Bold = $h1.Font.Bold
$headerTarget.HorizontalAlignment = $h1.HorizontalAlignment
$headerTarget.VerticalAlignment = $h1.VerticalAlignment
$headerTarget.Borders.LineStyle = 1

# Data values for columns I (I0) and J (IF), rows 2..59
$iVals = @(6,9,7,7,6,6,5,7,7,7,6,7,6,6,6,7,8,5,5,4,5,8,8,5,9,5,9,6,8,7,7,8,6,8,6,9,8,5,7,9,7,6,7,4,7,8,6,6,5,6,6,8,7,9,6,6,9,7)
$jVals = @(7,9,7,7,6,7,5,7,7,7,6,7,6,7,7,7,8,5,6,5,5,9,8,5,9,6,9,6,8,7,7,8,6,8,6,9,8,5,7,9,7,7,7,5,7,9,6,6,6,6,6,8,7,9,6,6,9,7)

for ($r = 2; $r -le 59; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
